$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1549
$ws.Range("F3").Value = 914
$ws.Range("F4").Value = 473
$ws.Range("F5").Value = 915
$ws.Range("F6").Value = 523
$ws.Range("F7").Value = 7838
$ws.Range("F11").Value = 5691
$ws.Range("F12").Value = 579
$ws.Range("F14").Value = 7917
$ws.Range("F15").Value = 9307
$ws.Range("F17").Value = 924
$ws.Range("F18").Value = 4541
$ws.Range("F19").Value = 688
$ws.Range("F20").Value = 261
$ws.Range("F21").Value = 85
$ws.Range("F24").Value = 1212
$ws.Range("F26").Value = 1704
$ws.Range("F27").Value = 744
$ws.Range("F28").Value = 965
$ws.Range("F29").Value = 15
$ws.Range("F30").Value = 1902
$ws.Range("F32").Value = 2349
$ws.Range("F34").Value = 1498
$ws.Range("F39").Value = 523
$ws.Range("F40").Value = 3014
$ws.Range("F41").Value = 4167
$ws.Range("F43").Value = 52
$ws.Range("F44").Value = 434
$ws.Range("F45").Value = 516
$ws.Range("F47").Value = 863
$ws.Range("F48").Value = 182
$ws.Range("F49").Value = 4108

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 17
$ws.Range("F15").Value = 50
$ws.Range("F25").Value = 62

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5346

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1549
$ws.Range("F4").Value = 914
$ws.Range("F5").Value = 473
$ws.Range("F6").Value = 915
$ws.Range("F7").Value = 523
$ws.Range("F10").Value = 5691
$ws.Range("F11").Value = 579
$ws.Range("F12").Value = 7917
$ws.Range("F13").Value = 9307
$ws.Range("F15").Value = 924
$ws.Range("F16").Value = 4541
$ws.Range("F17").Value = 688
$ws.Range("F18").Value = 261
$ws.Range("F19").Value = 85
$ws.Range("F23").Value = 1212
$ws.Range("F25").Value = 1704
$ws.Range("F26").Value = 744
$ws.Range("F27").Value = 965
$ws.Range("F28").Value = 15
$ws.Range("F29").Value = 1902
$ws.Range("F31").Value = 2349
$ws.Range("F37").Value = 62
$ws.Range("F39").Value = 523
$ws.Range("F40").Value = 4167
$ws.Range("F43").Value = 52
$ws.Range("F44").Value = 434
$ws.Range("F45").Value = 516
$ws.Range("F47").Value = 863
$ws.Range("F48").Value = 182
$ws.Range("F49").Value = 4108
